$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")
$newRange = $ws.Range("A8:K139")
$lo.Resize($newRange)
Write-Output $lo.Range.Address()
